$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 6 de Julio de 2020 a las 17:07"

# Kenia overtook El Salvador in total cases, so the two countries swap rows
# (the sheet is kept sorted by total cases, column B, descending).
$ws.Range("A76").Value = "Kenia"
$ws.Range("A77").Value = "El Salvador"

# Updated per-country statistics (Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes)
# row 4
$ws.Range("B4").Value = 2992605
$ws.Range("C4").Value = 9677
$ws.Range("D4").Value = 1290682
$ws.Range("E4").Value = 1569261
$ws.Range("G4").Value = 93
$ws.Range("H4").Value = 132662
# row 6
$ws.Range("B6").Value = 712920
$ws.Range("C6").Value = 15084
$ws.Range("D6").Value = 435441
$ws.Range("E6").Value = 257406
$ws.Range("G6").Value = 373
$ws.Range("H6").Value = 20073
# row 16
$ws.Range("B16").Value = 213716
$ws.Range("C16").Value = 4207
$ws.Range("D16").Value = 149634
$ws.Range("E16").Value = 62114
$ws.Range("G16").Value = 52
$ws.Range("H16").Value = 1968
# row 18
$ws.Range("B18").Value = 197651
$ws.Range("C18").Value = 93
$ws.Range("E18").Value = 6365
# row 57
$ws.Range("B57").Value = 20837
$ws.Range("C57").Value = 513
$ws.Range("D57").Value = 12182
$ws.Range("E57").Value = 8397
$ws.Range("G57").Value = 8
$ws.Range("H57").Value = 258
# row 61
$ws.Range("D61").Value = 11047
$ws.Range("E61").Value = 6179
$ws.Range("G61").Value = 3
$ws.Range("H61").Value = 588
# row 76
$ws.Range("B76").Value = 8067
$ws.Range("C76").Value = 181
$ws.Range("D76").Value = 2414
$ws.Range("E76").Value = 5489
$ws.Range("G76").Value = 4
$ws.Range("H76").Value = 164
# row 77
$ws.Range("B77").Value = 8027
$ws.Range("C77").Value = 250
$ws.Range("D77").Value = 4730
$ws.Range("E77").Value = 3074
$ws.Range("G77").Value = 6
$ws.Range("H77").Value = 223
# row 96
$ws.Range("B96").Value = 4341
$ws.Range("C96").Value = 64
$ws.Range("E96").Value = 3833
# row 112
$ws.Range("B112").Value = 2077
$ws.Range("C112").Value = 1
$ws.Range("E112").Value = 149
# row 115
$ws.Range("B115").Value = 1885
$ws.Range("C115").Value = 12
$ws.Range("E115").Value = 538
# row 139
$ws.Range("D139").Value = 270
$ws.Range("E139").Value = 709
# row 162
$ws.Range("B162").Value = 342
$ws.Range("C162").Value = 1
$ws.Range("E162").Value = 2
# row 164
$ws.Range("B164").Value = 316
$ws.Range("C164").Value = 3
$ws.Range("D164").Value = 245
$ws.Range("E164").Value = 65
